# update on 20210731 画中人
# Replace curly/smart quotes used around nicknames/titles with straight
# single quotes (and one straight double-quote pair) in several dialogue
# lines on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rsquo = [char]0x2019  # ’
$em    = [char]0x2014  # —

$ws.Range("C17").Value = "[name=`"Greatmouth Mob`"]  This is the most destructive match the Roar Arena has seen this event$em! The Butcher of the Field! 'Brassrust' Ingra!  `n"

$ws.Range("C18").Value = "[name=`"Greatmouth Mob`"]  Just last month, at the season${rsquo}s first event, poor 'Foehn' got absolutely decimated when the 'Brassrust' kept pummeling him until all of his limbs were crushed!`n"

$ws.Range("C21").Value = "[name=`"Greatmouth Mob`"]  Just like that! 'Brassrust' Ingra! Smear the battlefield with your opponent${rsquo}s blood!  `n"

$ws.Range("C62").Value = " 'To be a knight is to be the noble light that illuminates the land...' `n"

$ws.Range("C116").Value = "[name=`"Platinum`"]  ...This one has a pretty name. Let${rsquo}s go with this one, the 'Thorn Tear.'`n"
